# Update the "SW" sheet with two additional blocks of correlation stats
# (margin-zone and inland-zone), mirroring the existing bare/dark-ice block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SW")

# ---- Row 1: group headers (margin / inland) + bare ice / dark ice sub-headers ----
$ws.Range("H1").Value = "margin"
$ws.Range("I1").Value = "bare ice"
$ws.Range("J1").Value = ""
$ws.Range("K1").Value = "dark ice"
$ws.Range("L1").Value = ""

$ws.Range("N1").Value = "inland"
$ws.Range("O1").Value = "bare ice"
$ws.Range("P1").Value = ""
$ws.Range("Q1").Value = "dark ice"
$ws.Range("R1").Value = ""

# ---- Row 2: R / p sub headers ----
$ws.Range("I2").Value = "R"
$ws.Range("J2").Value = "p"
$ws.Range("K2").Value = "R"
$ws.Range("L2").Value = "p"

$ws.Range("O2").Value = "R"
$ws.Range("P2").Value = "p"
$ws.Range("Q2").Value = "R"
$ws.Range("R2").Value = "p"

# ---- Row labels (elevation, slope, aspect, distance, duration) ----
$ws.Range("H3").Value = "elevation"
$ws.Range("H4").Value = "slope"
$ws.Range("H5").Value = "aspect"
$ws.Range("H6").Value = "distance"
$ws.Range("H7").Value = "duration"

$ws.Range("N3").Value = "elevation"
$ws.Range("N4").Value = "slope"
$ws.Range("N5").Value = "aspect"
$ws.Range("N6").Value = "distance"
$ws.Range("N7").Value = "duration"

# ---- margin block values (bare ice R/p, dark ice R/p) ----
$ws.Range("I3").Value = 0.070000000000000007
$ws.Range("J3").Value = 0.88
$ws.Range("K3").Value = -0.77
$ws.Range("L3").Value = 0.04

$ws.Range("I4").Value = 0.3
$ws.Range("J4").Value = 0.51
$ws.Range("K4").Value = 0.82
$ws.Range("L4").Value = 0.02

$ws.Range("I5").Value = -0.28999999999999998
$ws.Range("J5").Value = 0.53
$ws.Range("K5").Value = -0.67
$ws.Range("L5").Value = 0.1

$ws.Range("I6").Value = -0.63
$ws.Range("J6").Value = 0.13
$ws.Range("K6").Value = -0.56000000000000005
$ws.Range("L6").Value = 0.19

$ws.Range("I7").Value = 0.69
$ws.Range("J7").Value = 0.09
$ws.Range("K7").Value = -0.49
$ws.Range("L7").Value = 0.26

# ---- inland block values (bare ice R/p, dark ice R/p) ----
$ws.Range("O3").Value = 0.27
$ws.Range("P3").Value = 0.55000000000000004
$ws.Range("Q3").Value = -0.17
$ws.Range("R3").Value = 0.71

$ws.Range("O4").Value = 0.08
$ws.Range("P4").Value = 0.86
$ws.Range("Q4").Value = 0.04
$ws.Range("R4").Value = 0.94

$ws.Range("O5").Value = -0.56999999999999995
$ws.Range("P5").Value = 0.18
$ws.Range("Q5").Value = 0.25
$ws.Range("R5").Value = 0.6

$ws.Range("O6").Value = 0.45
$ws.Range("P6").Value = 0.31
$ws.Range("Q6").Value = 0.09
$ws.Range("R6").Value = 0.85

$ws.Range("O7").Value = 0.01
$ws.Range("P7").Value = 0.99
$ws.Range("Q7").Value = 0.03
$ws.Range("R7").Value = 0.96

# ---- number formatting for the new R / p columns (match the existing columns) ----
$ws.Range("I3:J7").NumberFormat = "0.00"
$ws.Range("K3:L7").NumberFormat = "0.00"
$ws.Range("O3:P7").NumberFormat = "0.00"
$ws.Range("Q3:R7").NumberFormat = "0.00"

$ws.Range("I3:J7").Font.Bold = $false
$ws.Range("K3:L7").Font.Bold = $false

# ---- alignment: center the group/sub headers (also re-applies to B1:E1 so it
#      matches the refreshed header style used elsewhere in the workbook) ----
$ws.Range("B1:E1").HorizontalAlignment = -4108
$ws.Range("B1:E1").VerticalAlignment = -4108
$ws.Range("H1:L1").HorizontalAlignment = -4108
$ws.Range("H1:L1").VerticalAlignment = -4108
$ws.Range("N1:R1").HorizontalAlignment = -4108
$ws.Range("N1:R1").VerticalAlignment = -4108

$ws.Range("I2:L2").HorizontalAlignment = -4108
$ws.Range("I2:L2").VerticalAlignment = -4108
$ws.Range("O2:R2").HorizontalAlignment = -4108
$ws.Range("O2:R2").VerticalAlignment = -4108
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("N2").HorizontalAlignment = -4108
$ws.Range("N2").VerticalAlignment = -4108

$ws.Range("H3:H7").HorizontalAlignment = -4108
$ws.Range("H3:H7").VerticalAlignment = -4108
$ws.Range("N3:N7").HorizontalAlignment = -4108
$ws.Range("N3:N7").VerticalAlignment = -4108

# ---- merge the two-column group headers ----
$ws.Range("I1:J1").Merge()
$ws.Range("K1:L1").Merge()
$ws.Range("O1:P1").Merge()
$ws.Range("Q1:R1").Merge()

# ---- column / view bits so the sheet lands on the new data ----
$ws.Columns.Item("C").AutoFit() | Out-Null
$ws.Range("Q8").Select() | Out-Null

Write-Host "SW sheet updated with margin/inland correlation blocks"
